# "multiple row write in excel updated"
# Adds a new data row (row 5) to the GBT test-data sheet, mirroring the
# layout of the existing row 2 (RuleType/ClassName/RuleName/RuleSetVersion
# in A:D, a Param value in E, an Expected/Actual/Pass-Fail triple in Q:S),
# and leaves the selection on A4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text columns for the new row.
$ws.Range("A5").Value = "Decision_Table"
$ws.Range("B5").Value = "PegaFS-Data-RelCodes"
$ws.Range("C5").Value = "RelatedPartyEnforcedPairs"
$ws.Range("D5").Value = "PegaFS:08-06-01"
$ws.Range("E5").Value = "aff"

# Q5 ("Expected" column) should hold the text "true" formatted exactly like
# Q2 (text-as-string with a quote-prefix/Text number format). Copying Q2's
# value+format onto Q5 reproduces that precisely.
$ws.Range("Q2").Copy($ws.Range("Q5"))

# R5 ("Actual" column) should hold the literal text "false" (not the
# boolean FALSE). Assigning the string directly gets auto-coerced to a
# boolean cell, so instead we compute it as a formula result (a text
# value) and then flatten the formula down to a plain value via copy /
# paste-values, which keeps it as text.
$ws.Range("R5").Formula = "=""false"""
$ws.Range("R5").Copy()
$ws.Range("R5").PasteSpecial(-4163)  # xlPasteValues

# S5 ("Pass/Fail" column).
$ws.Range("S5").Value = "Fail"

# Leave the active selection on A4, as in the edited workbook.
[void]$ws.Range("A4").Select()
